$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append a "?" to the end of the "Boxplot fuer Interaktionszeiten"
#    paragraph's text (same run formatting as the rest of the
#    paragraph).
# ------------------------------------------------------------------
$boxplotIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "Boxplot fuer Interaktionszeiten`r") {
        $boxplotIndex = $i
        break
    }
}

$pBoxplot = $d.Paragraphs.Item($boxplotIndex)
$rBoxplot = $pBoxplot.Range
$rBoxplot.End = $rBoxplot.End - 1
$boxplotTextEnd = $rBoxplot.End
$d.Range($boxplotTextEnd, $boxplotTextEnd).InsertAfter("?")

# ------------------------------------------------------------------
# 2) Insert a brand-new paragraph right after it and give it the
#    text "Boxplot für Vorerfahrung?" (formatting is inherited from
#    the paragraph we split off from).
# ------------------------------------------------------------------
$pBoxplot = $d.Paragraphs.Item($boxplotIndex)
$pBoxplot.Range.InsertParagraphAfter()

$pNew = $d.Paragraphs.Item($boxplotIndex + 1)
$pNew.Range.Text = "Boxplot für Vorerfahrung?"

# ------------------------------------------------------------------
# 3) Move the (hidden) "_GoBack" bookmark from wherever it currently
#    sits onto the end of the freshly typed text in the new
#    paragraph, i.e. right before its paragraph mark.
#
#    NOTE: creating a zero-length bookmark exactly at
#    "paragraph.End - 1" (the paragraph-mark offset) confuses the
#    engine, so we temporarily insert a throw-away character after
#    the target spot, anchor the bookmark there (now a safe, non
#    boundary offset), and then delete the throw-away character
#    again. This leaves the bookmark collapsed at the correct spot.
# ------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$rNew = $pNew.Range
$rNew.End = $rNew.End - 1
$targetPos = $rNew.End
$d.Range($targetPos, $targetPos).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($targetPos, $targetPos))
$d.Range($targetPos, $targetPos + 1).Text = ""

Write-Output "done"
